# Apply scheduled runner updates to Exodus_Profits sheets
$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3132.3333
$ws.Range("J100").Value = 2758.8
$ws.Range("L100").Value = 2758.8
$ws.Range("N100").Value = -3840.8

# Hunk 1: ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 372905.28
$ws.Range("I106").Value = 696560.4
$ws.Range("K106").Value = 696560.4
$ws.Range("M106").Value = -695929.4

# Hunk 2: ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 92314.46000000001
$ws.Range("J133").Value = 92314.46000000001
$ws.Range("L133").Value = 92314.46000000001
$ws.Range("N133").Value = -102434.46

# Hunk 3: ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 57885.715
$ws.Range("J134").Value = 57885.715
$ws.Range("L134").Value = 57885.715
$ws.Range("N134").Value = -68025.715

# Hunk 4: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 605.5417
$ws.Range("I135").Value = 666.4737
$ws.Range("J135").Value = 374
$ws.Range("K135").Value = 5998.263300000001
$ws.Range("L135").Value = 3366
$ws.Range("M135").Value = -3463.263300000001
$ws.Range("N135").Value = -8436

# Hunk 5: ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 78831.71000000001
$ws.Range("J136").Value = 78831.71000000001
$ws.Range("L136").Value = 78831.71000000001
$ws.Range("N136").Value = -89031.71000000001

# Hunk 6: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 632205.8
$ws.Range("I137").Value = 1983.9412
$ws.Range("J137").Value = 2417834.2
$ws.Range("K137").Value = 5951.8236
$ws.Range("L137").Value = 7253502.600000001
$ws.Range("M137").Value = -3401.8236
$ws.Range("N137").Value = -7258602.600000001

# Hunk 7: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3090.3809
$ws.Range("I138").Value = 1775.5
$ws.Range("K138").Value = 5326.5
$ws.Range("M138").Value = -186.5

# Hunk 8: ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270

# Hunk 9: ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 80431.5
$ws.Range("J140").Value = 80778.86
$ws.Range("L140").Value = 80778.86
$ws.Range("N140").Value = -91138.86

# Hunk 10: ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1038.4445
$ws.Range("I2").Value = 857.4545000000001
$ws.Range("J2").Value = 1322.8572
$ws.Range("K2").Value = 857.4545000000001
$ws.Range("L2").Value = 1322.8572
$ws.Range("M2").Value = -744.4545000000001
$ws.Range("N2").Value = -1548.8572

# Hunk 11: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24220.37
$ws.Range("I32").Value = 13666.667
$ws.Range("K32").Value = 13666.667
$ws.Range("M32").Value = -13379.667

# Hunk 12: ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2996.3635
$ws.Range("I45").Value = 2687.1428
$ws.Range("K45").Value = 2687.1428
$ws.Range("M45").Value = -2310.1428

# Hunk 13: ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1038.4445
$ws.Range("I116").Value = 857.4545000000001
$ws.Range("J116").Value = 1322.8572
$ws.Range("K116").Value = 857.4545000000001
$ws.Range("L116").Value = 1322.8572
$ws.Range("M116").Value = 1436.5455
$ws.Range("N116").Value = -5910.8572

# Hunk 14: ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 89954.75
$ws.Range("J139").Value = 89954.75
$ws.Range("L139").Value = 89954.75
$ws.Range("N139").Value = -100234.75

# Hunk 15: BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1038.4445
$ws.Range("I3").Value = 857.4545000000001
$ws.Range("J3").Value = 1322.8572
$ws.Range("K3").Value = 857.4545000000001
$ws.Range("L3").Value = 1322.8572
$ws.Range("M3").Value = -743.4545000000001
$ws.Range("N3").Value = -1550.8572

# Hunk 16: BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1468886.8
$ws.Range("I99").Value = 61626.94
$ws.Range("K99").Value = 61626.94
$ws.Range("M99").Value = -60128.94

# Hunk 17: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4969.4287
$ws.Range("I107").Value = 4646.2
$ws.Range("J107").Value = 5777.5
$ws.Range("K107").Value = 4646.2
$ws.Range("L107").Value = 5777.5
$ws.Range("M107").Value = -2726.2
$ws.Range("N107").Value = -9617.5

# Hunk 18: BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 60390
$ws.Range("J130").Value = 60390
$ws.Range("L130").Value = 60390
$ws.Range("N130").Value = -70430

# Hunk 19: BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 81272.5
$ws.Range("J132").Value = 81272.5
$ws.Range("L132").Value = 81272.5
$ws.Range("N132").Value = -91392.5

# Hunk 20: BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 88781.11
$ws.Range("J135").Value = 88781.11
$ws.Range("L135").Value = 88781.11
$ws.Range("N135").Value = -98921.11

# Hunk 21: BSM row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 79897.664
$ws.Range("J137").Value = 79897.664
$ws.Range("L137").Value = 79897.664
$ws.Range("N137").Value = -90097.664

# Hunk 22: BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 86558.11
$ws.Range("J138").Value = 86558.11
$ws.Range("L138").Value = 86558.11
$ws.Range("N138").Value = -96838.11

# Hunk 23: BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 100355.336
$ws.Range("J140").Value = 64024
$ws.Range("L140").Value = 64024
$ws.Range("N140").Value = -74384

# Hunk 24: CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 36334.38
$ws.Range("I134").Value = 1702.8695
$ws.Range("K134").Value = 5108.6085
$ws.Range("M134").Value = -2573.6085

# Hunk 25: CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 99814.28999999999
$ws.Range("J138").Value = 99814.28999999999
$ws.Range("L138").Value = 99814.28999999999
$ws.Range("N138").Value = -110094.29

# Hunk 26: CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5336.3335
$ws.Range("I56").Value = 5336.3335
$ws.Range("K56").Value = 5336.3335
$ws.Range("M56").Value = -4806.3335

# Hunk 27: CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 493.5
$ws.Range("I107").Value = 493.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1480.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 439.5
$ws.Range("N107").ClearContents()

# Hunk 28: CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 6796.8184
$ws.Range("J114").Value = 7376.5
$ws.Range("L114").Value = 22129.5
$ws.Range("N114").Value = -28637.5

# Hunk 29: CUL row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 2382.5715
$ws.Range("J127").Value = 2382.5715
$ws.Range("L127").Value = 7147.7145
$ws.Range("N127").Value = -17067.7145

# Hunk 30: GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2460.3684
$ws.Range("I97").Value = 859.25
$ws.Range("K97").Value = 859.25
$ws.Range("M97").Value = -363.25

# Hunk 31: GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 461.05884
$ws.Range("I107").Value = 288.5
$ws.Range("K107").Value = 288.5
$ws.Range("M107").Value = 1631.5

# Hunk 32: GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3131.96
$ws.Range("I126").Value = 2749.5
$ws.Range("K126").Value = 8248.5
$ws.Range("M126").Value = -5778.5

# Hunk 33: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4294.1665
$ws.Range("I132").Value = 4849.8335
$ws.Range("K132").Value = 14549.5005
$ws.Range("M132").Value = -12019.5005

# Hunk 34: GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 96658.336
$ws.Range("J135").Value = 96658.336
$ws.Range("L135").Value = 96658.336
$ws.Range("N135").Value = -106798.336

# Hunk 35: GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 64446.668
$ws.Range("J140").Value = 66752.5
$ws.Range("L140").Value = 66752.5
$ws.Range("N140").Value = -77112.5

# Hunk 36: LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1599.6
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1599.6
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1599.6
$ws.Range("N46").Value = -1975.6
$ws.Range("M46").ClearContents()

# Hunk 37: LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 29937.375
$ws.Range("I100").Value = 38249.832
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 38249.832
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -37708.832
$ws.Range("N100").Value = -6082

# Hunk 38: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3184.0588
$ws.Range("I136").Value = 2821.9333
$ws.Range("J136").Value = 3469.9473
$ws.Range("K136").Value = 8465.7999
$ws.Range("L136").Value = 10409.8419
$ws.Range("M136").Value = -5915.7999
$ws.Range("N136").Value = -15509.8419

# Hunk 39: WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1631.6571
$ws.Range("I107").Value = 961.3158
$ws.Range("J107").Value = 2427.6875
$ws.Range("K107").Value = 2883.9474
$ws.Range("L107").Value = 7283.0625
$ws.Range("M107").Value = -963.9474
$ws.Range("N107").Value = -11123.0625

# Hunk 40: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1812765.4
$ws.Range("I132").Value = 669.3125
$ws.Range("J132").Value = 5436957.5
$ws.Range("K132").Value = 2007.9375
$ws.Range("L132").Value = 16310872.5
$ws.Range("M132").Value = 522.0625
$ws.Range("N132").Value = -16315932.5

# Hunk 41: WVR row 140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 84816.5
$ws.Range("J140").Value = 86379.8
$ws.Range("L140").Value = 86379.8
$ws.Range("N140").Value = -96739.8
